$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-07-29"

# Update the header label (shared string used by cell I1)
$ws.Range("I1").Value = "2022 (through 07-29)"

# Update the updated data values
$ws.Range("I8").Value = 162
$ws.Range("I14").Value = 968
